# "The Cake is a Lie.xlsx" - add SPARC condition-code ("cc") and extended (X)
# opcode rows to the operation table on Hoja1.
#
# The table originally listed 8 base opcodes (rows 3-10). The edit expands it
# to 20 opcodes (rows 3-22), inserting the corresponding *cc / X / Xcc variants
# next to their base instruction, re-sorts a couple of rows, normalizes every
# data row to the "bordered + centered + text-format" style, and narrows the
# active selection back down to the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Normalize formatting for the whole data block first (copy/paste formats
#    only, so no cell value / shared-string is touched here).
#    - C3:E14  -> bordered, centered horiz+vert, text number format (style 2)
#    - C15:D22 -> same style as above (style 2)
#    - E15:E22 -> bordered, centered horiz+vert, general format (style 1)
# ---------------------------------------------------------------------------
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C3:E14").PasteSpecial(-4122) | Out-Null

$ws.Range("C15:D22").PasteSpecial(-4122) | Out-Null

$ws.Range("E2").Copy() | Out-Null
$ws.Range("E15:E22").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Fill in the op / op3 / operation columns for rows 3-22.
#    Values that already existed in the shared-string table are written
#    first; brand-new values are written afterwards in the same order they
#    were introduced so the rebuilt table matches the original edit.
# ---------------------------------------------------------------------------

# -- values that reuse already-existing shared strings --
$ws.Range("C3").Value = "01"
$ws.Range("D3").Value = "000001"
$ws.Range("E3").Value = "AND"
$ws.Range("C4").Value = "01"
$ws.Range("C5").Value = "01"
$ws.Range("D5").Value = "000101"
$ws.Range("E5").Value = "ANDN"
$ws.Range("C6").Value = "01"
$ws.Range("C7").Value = "01"
$ws.Range("D7").Value = "000010"
$ws.Range("E7").Value = "OR"
$ws.Range("C8").Value = "01"
$ws.Range("C9").Value = "01"
$ws.Range("D9").Value = "000110"
$ws.Range("E9").Value = "ORN"
$ws.Range("C10").Value = "01"
$ws.Range("C11").Value = "01"
$ws.Range("D11").Value = "000011"
$ws.Range("E11").Value = "XOR"
$ws.Range("C12").Value = "01"
$ws.Range("C13").Value = "01"
$ws.Range("D13").Value = "000111"
$ws.Range("E13").Value = "XNOR"
$ws.Range("C14").Value = "01"
$ws.Range("C15").Value = "01"
$ws.Range("D15").Value = "000000"
$ws.Range("E15").Value = "ADD"
$ws.Range("C16").Value = "01"
$ws.Range("C17").Value = "01"
$ws.Range("C18").Value = "01"
$ws.Range("C19").Value = "01"
$ws.Range("D19").Value = "000100"
$ws.Range("E19").Value = "SUB"
$ws.Range("C20").Value = "01"
$ws.Range("C21").Value = "01"
$ws.Range("C22").Value = "01"

# -- brand-new values (new unique shared strings) --
$ws.Range("E4").Value = "ANDcc"
$ws.Range("D6").Value = "010101"
$ws.Range("E6").Value = "ANDNcc"
$ws.Range("D4").Value = "010001"
$ws.Range("D10").Value = "010110"
$ws.Range("E10").Value = "ORNcc"
$ws.Range("D8").Value = "010010"
$ws.Range("E8").Value = "ORcc"
$ws.Range("D14").Value = "010111"
$ws.Range("E14").Value = "XNORcc"
$ws.Range("D12").Value = "010011"
$ws.Range("E12").Value = "XORcc"
$ws.Range("D18").Value = "011000"
$ws.Range("E18").Value = "ADDXcc"
$ws.Range("D17").Value = "001000"
$ws.Range("E17").Value = "ADDX"
$ws.Range("D16").Value = "010000"
$ws.Range("E16").Value = "ADDcc"
$ws.Range("D22").Value = "011100"
$ws.Range("E22").Value = "SUBXcc"
$ws.Range("D21").Value = "001100"
$ws.Range("E21").Value = "SUBX"
$ws.Range("D20").Value = "010100"
$ws.Range("E20").Value = "SUBcc"

# ---------------------------------------------------------------------------
# 3) Fix up the selection: it now only covers the header row, not the
#    whole (now much larger) table.
# ---------------------------------------------------------------------------
$ws.Range("C2:E2").Select() | Out-Null
